$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scenario")

# --- Update existing Execution Status values: "Yes" -> "No" for Checkout & Contact Us blocks ---
$ws.Range("F21").Value = "No"
$ws.Range("F22").Value = "No"
$ws.Range("F23").Value = "No"
$ws.Range("F25").Value = "No"
$ws.Range("F26").Value = "No"

# --- Append new "Book Flights" scenario block in rows 28-30 ---
# Copy formatting from the analogous existing block (rows 25-27: a single-row
# header, a detail row with wrapped description, then a blank separator row)
# so the new rows reuse the workbook's existing cell styles instead of
# creating new ones.
$ws.Range("A25:F25").Copy()
$ws.Range("A28:F28").PasteSpecial(-4122)

$ws.Range("A26:F26").Copy()
$ws.Range("A29:F29").PasteSpecial(-4122)

$ws.Range("A27:F27").Copy()
$ws.Range("A30:F30").PasteSpecial(-4122)

$excel.CutCopyMode = 0

$ws.Range("A28").Value = "Book Flights"
$ws.Range("F28").Value = "Yes"

$ws.Range("B29").Value = "Select_Destination"
$ws.Range("C29").Value = "Book Flights"
$ws.Range("D29").Value = "Test1Flight"
$ws.Range("E29").Value = "Verify booking flights successfully"
$ws.Range("F29").Value = "Yes"

# Row 30 stays blank (separator row), matching the pattern used throughout the sheet.

# --- Update the view to match the recorded selection/scroll position ---
$ws.Range("C33").Select()
$excel.ActiveWindow.ScrollRow = 19
